$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 208 (1-based), pushing the existing
# rows 208..288 down to 209..289.
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with its data.
$ws.Range("A208").Value = 9
$ws.Range("B208").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C208").Value = "Metropolitana"
$ws.Range("D208").Value = 44452
$ws.Range("E208").Value = 13
$ws.Range("F208").Value = 100112024
$ws.Range("G208").Value = "Choclo"
$ws.Range("H208").Value = "Dulce o Americano"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 52
$ws.Range("K208").Value = 39000
$ws.Range("L208").Value = 40000
$ws.Range("M208").Value = 39500
$ws.Range("N208").Value = "$/malla 70 unidades"
$ws.Range("O208").Value = "Región de Arica y Parinacota"
$ws.Range("P208").Value = 564
$ws.Range("Q208").Value = 70
$ws.Range("R208").Value = "Hortaliza"
